$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename event "NEON BOUNCE" -> "STROBO BOUNCE" (cell B359)
# ---------------------------------------------------------------------------
$ws.Range("B359").Value = "STROBO BOUNCE"

# ---------------------------------------------------------------------------
# Helper data for the six new rows (462-467)
# Columns: A=Datum (serial date), B=Event, C=Location, D=Stadt, E=Link(url)
# ---------------------------------------------------------------------------

# ---- Row 462 : RAVE BOAT --------------------------------------------------
$ws.Range("A461").Copy()
$ws.Range("A462").PasteSpecial(-4122)
$ws.Range("A462").Value = 45815

$ws.Range("B461").Copy()
$ws.Range("B462").PasteSpecial(-4122)
$ws.Range("B462").Value = "RAVE BOAT"

$ws.Range("C461").Copy()
$ws.Range("C462").PasteSpecial(-4122)
$ws.Range("C462").Value = "Marina Duisburg"

$ws.Range("D461").Copy()
$ws.Range("D462").PasteSpecial(-4122)
$ws.Range("D462").Value = "Duisburg"

$url462 = "https://stagedates.com/events/boat-rave-marina-duisburg-marina-duisburg-20250607-CKeQl"
$ws.Hyperlinks.Add($ws.Range("E462"), $url462, "", "", $url462)
$ws.Range("B461").Copy()
$ws.Range("E462").PasteSpecial(-4122)
$ws.Range("E462").Value = $url462
$len462 = $url462.Length
$c1 = $ws.Range("E462").Characters(1, $len462 - 1)
$c1.Font.Underline = 2
$c1.Font.ColorIndex = 4
$c1.Font.Name = "Calibri"
$c1.Font.Size = 11
$c2 = $ws.Range("E462").Characters($len462, 1)
$c2.Font.Underline = 2
$c2.Font.ColorIndex = 4
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11

# ---- Row 463 : POLAAR OPEN AIR --------------------------------------------
$ws.Range("A461").Copy()
$ws.Range("A463").PasteSpecial(-4122)
$ws.Range("A463").Value = 45899

$ws.Range("B461").Copy()
$ws.Range("B463").PasteSpecial(-4122)
$ws.Range("B463").Value = "POLAAR OPEN AIR"

$ws.Range("C461").Copy()
$ws.Range("C463").PasteSpecial(-4122)
$ws.Range("C463").Value = "Haus Witten"

$ws.Range("D461").Copy()
$ws.Range("D463").PasteSpecial(-4122)
$ws.Range("D463").Value = "Witten"

$url463 = "https://polaar.ticket.io/brfmsfz7/"
$ws.Hyperlinks.Add($ws.Range("E463"), $url463, "", "", $url463)
$ws.Range("B461").Copy()
$ws.Range("E463").PasteSpecial(-4122)
$ws.Range("E463").Value = $url463
$len463 = $url463.Length
$c1 = $ws.Range("E463").Characters(1, $len463 - 1)
$c1.Font.Underline = 2
$c1.Font.ColorIndex = 4
$c1.Font.Name = "Calibri"
$c1.Font.Size = 11
$c2 = $ws.Range("E463").Characters($len463, 1)
$c2.Font.Underline = 2
$c2.Font.ColorIndex = 4
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11

# ---- Row 464 : 222' BERGFEST CPDO -----------------------------------------
$ws.Range("A461").Copy()
$ws.Range("A464").PasteSpecial(-4122)
$ws.Range("A464").Value = 45798

$ws.Range("B461").Copy()
$ws.Range("B464").PasteSpecial(-4122)
$ws.Range("B464").Value = "222‘ BERGFEST CPDO"

$ws.Range("C461").Copy()
$ws.Range("C464").PasteSpecial(-4122)
$ws.Range("C464").Value = "Prismatic"

$ws.Range("D461").Copy()
$ws.Range("D464").PasteSpecial(-4122)
$ws.Range("D464").Value = "Dortmund"

$url464 = "https://www.instagram.com/reel/DJth1tftGST/?igsh=Z3J3ZWZ3OGI1OHcz"
$ws.Hyperlinks.Add($ws.Range("E464"), $url464, "", "", $url464)
$ws.Range("B461").Copy()
$ws.Range("E464").PasteSpecial(-4122)
$ws.Range("E464").Value = $url464
$len464 = $url464.Length
$c1 = $ws.Range("E464").Characters(1, $len464 - 1)
$c1.Font.Underline = 2
$c1.Font.ColorIndex = 4
$c1.Font.Name = "Calibri"
$c1.Font.Size = 11
$c2 = $ws.Range("E464").Characters($len464, 1)
$c2.Font.Underline = 2
$c2.Font.ColorIndex = 4
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11

# ---- Row 465 : #MITTWOCHENENDE --------------------------------------------
$ws.Range("A461").Copy()
$ws.Range("A465").PasteSpecial(-4122)
$ws.Range("A465").Value = 45798

$ws.Range("B461").Copy()
$ws.Range("B465").PasteSpecial(-4122)
$ws.Range("B465").Value = "#MITTWOCHENENDE"

$ws.Range("C461").Copy()
$ws.Range("C465").PasteSpecial(-4122)
$ws.Range("C465").Value = "Odonien"

$ws.Range("D461").Copy()
$ws.Range("D465").PasteSpecial(-4122)
$ws.Range("D465").Value = "Köln"

$url465 = "https://www.instagram.com/p/DJru0uiswOw/?igsh=MXhoMThnMm42NjNhZQ=="
$ws.Hyperlinks.Add($ws.Range("E465"), $url465, "", "", $url465)
$ws.Range("B461").Copy()
$ws.Range("E465").PasteSpecial(-4122)
$ws.Range("E465").Value = $url465
$len465 = $url465.Length
$c1 = $ws.Range("E465").Characters(1, $len465 - 1)
$c1.Font.Underline = 2
$c1.Font.ColorIndex = 4
$c1.Font.Name = "Calibri"
$c1.Font.Size = 11
$c2 = $ws.Range("E465").Characters($len465, 1)
$c2.Font.Underline = 2
$c2.Font.ColorIndex = 4
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11

# ---- Row 466 : BLOOD MOON ---------------------------------------------------
$ws.Range("A461").Copy()
$ws.Range("A466").PasteSpecial(-4122)
$ws.Range("A466").Value = 45805

$ws.Range("B461").Copy()
$ws.Range("B466").PasteSpecial(-4122)
$ws.Range("B466").Value = "BLOOD MOON"

$ws.Range("C461").Copy()
$ws.Range("C466").PasteSpecial(-4122)
$ws.Range("C466").Value = "Sam‘s"

$ws.Range("D461").Copy()
$ws.Range("D466").PasteSpecial(-4122)
$ws.Range("D466").Value = "Bielefeld"

$url466 = "https://www.instagram.com/reel/DJzmHXwu6iU/?igsh=eHl3Ync5M3pxcGwx"
$ws.Hyperlinks.Add($ws.Range("E466"), $url466, "", "", $url466)
$ws.Range("B461").Copy()
$ws.Range("E466").PasteSpecial(-4122)
$ws.Range("E466").Value = $url466
$len466 = $url466.Length
$c1 = $ws.Range("E466").Characters(1, $len466 - 1)
$c1.Font.Underline = 2
$c1.Font.ColorIndex = 4
$c1.Font.Name = "Calibri"
$c1.Font.Size = 11
$c2 = $ws.Range("E466").Characters($len466, 1)
$c2.Font.Underline = 2
$c2.Font.ColorIndex = 4
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11

# ---- Row 467 : BERRYLECTRO HARDTECHNO --------------------------------------
$ws.Range("A461").Copy()
$ws.Range("A467").PasteSpecial(-4122)
$ws.Range("A467").Value = 45857

$ws.Range("B461").Copy()
$ws.Range("B467").PasteSpecial(-4122)
$ws.Range("B467").Value = "BERRYLECTRO HARDTECHNO"

$ws.Range("C461").Copy()
$ws.Range("C467").PasteSpecial(-4122)
$ws.Range("C467").Value = "Favela"

$ws.Range("D461").Copy()
$ws.Range("D467").PasteSpecial(-4122)
$ws.Range("D467").Value = "Münster"

$url467 = "https://www.instagram.com/berrylectro?igsh=Ym0yeG5xaHFzemUx"
$ws.Hyperlinks.Add($ws.Range("E467"), $url467, "", "", $url467)
$ws.Range("B461").Copy()
$ws.Range("E467").PasteSpecial(-4122)
$ws.Range("E467").Value = $url467
$len467 = $url467.Length
$c1 = $ws.Range("E467").Characters(1, $len467 - 1)
$c1.Font.Underline = 2
$c1.Font.ColorIndex = 4
$c1.Font.Name = "Calibri"
$c1.Font.Size = 11
$c2 = $ws.Range("E467").Characters($len467, 1)
$c2.Font.Underline = 2
$c2.Font.ColorIndex = 4
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11
